# Apply "Added another notch filter lol" edits to RCValueSelection.xlsx
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New annotation in row 14 (G14): note about resistor value being too high ---
$ws.Range("G14").Value = "RESISTOR VALUE TOO HIGH, SINCE SCOPES ARE LIKE 10M OHMS, SO WEIRD WOBBLES"

# --- New annotation in row 25 (I25): note about mislabelled caps ---
$ws.Range("I25").Value = "THIS WAS THE KAKA THAT DIDN'T WORK BC MISLABELLED CAPS"

# --- Row 26: replace the old "try this value" testing note with the final result note ---
$ws.Range("F26").Value = "THIS WORKED, BUT NO 1% TOLERANCE RESISTORS AVAIL. AT THIS VALUE"

# --- Row 29: highlight the chosen R/C values in yellow + bold, drop the stale "10k" AKA ---
$ws.Range("A29").Font.Bold = $true
$ws.Range("A29").Interior.Color = 65535
$ws.Range("A29").NumberFormat = "0.00E+00"

$ws.Range("B29").Font.Bold = $true
$ws.Range("B29").Interior.Color = 65535
$ws.Range("B29").NumberFormat = "0.00E+00"

$ws.Range("C29").Font.Bold = $true
$ws.Range("C29").Interior.Color = 65535

$ws.Range("D29").ClearContents()
$ws.Range("D29").Font.Bold = $true

$ws.Range("E29").Font.Bold = $true

$ws.Range("F29").Value = "http://www.digikey.com/product-detail/en/stackpole-electronics-inc/RNF14FTD12K1/RNF14FTD12K1CT-ND/1975097"
$ws.Range("F29").Font.Bold = $true

# --- Row 30: second datasheet/shopping link for the resistor ---
$ws.Range("F30").Value = "http://www.digikey.com/products/en/resistors/chip-resistor-surface-mount/52?FV=40087%2Cc0001%2Cc0163%2Cc0165%2Cc0172%2Cc0179%2Cc017c%2Cc002a%2Cc01af%2Cffe00034&mnonly=0&newproducts=0&ColumnSort=1000008&page=1&stock=1&pbfree=1&rohs=1&quantity=&ptm=0&fid=0&pageSize=25"
$ws.Range("F30").Font.Bold = $true

# --- Match the editor's final cursor position ---
$ws.Range("A29").Select()
